# Update "想去人数" (number of people interested) figures on the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1028
$ws1.Range("F3").Value = 13530
$ws1.Range("F5").Value = 1028
$ws1.Range("F6").Value = 19
$ws1.Range("F7").Value = 1738
$ws1.Range("F8").Value = 139
$ws1.Range("F10").Value = 80
$ws1.Range("F11").Value = 41
$ws1.Range("F13").Value = 13545
$ws1.Range("F14").Value = 338
$ws1.Range("F16").Value = 8960
$ws1.Range("F18").Value = 8049
$ws1.Range("F20").Value = 10
$ws1.Range("F21").Value = 150
$ws1.Range("F23").Value = 147
$ws1.Range("F25").Value = 22
$ws1.Range("F30").Value = 207
$ws1.Range("F31").Value = 183

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 39

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1028
$ws4.Range("F3").Value = 13530
$ws4.Range("F5").Value = 1028
$ws4.Range("F6").Value = 19
$ws4.Range("F7").Value = 1738
$ws4.Range("F8").Value = 139
$ws4.Range("F10").Value = 80
$ws4.Range("F11").Value = 41
$ws4.Range("F13").Value = 13545
$ws4.Range("F14").Value = 338
$ws4.Range("F16").Value = 8960
$ws4.Range("F18").Value = 8049
$ws4.Range("F20").Value = 10
$ws4.Range("F21").Value = 150
$ws4.Range("F23").Value = 147
$ws4.Range("F25").Value = 22
$ws4.Range("F29").Value = 39
$ws4.Range("F32").Value = 207
$ws4.Range("F33").Value = 183
